$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Mark the run that holds the captured-waveform drawing as
#    "do not spell/grammar check" (<w:rPr><w:noProof/></w:rPr>).
# ------------------------------------------------------------------
$pic = $d.InlineShapes.Item(1)
$pic.Range.Font.NoProofing = -1

# ------------------------------------------------------------------
# 2) Append the new "May 1 2016" diary entry after the last
#    paragraph, preceded by one blank paragraph.
# ------------------------------------------------------------------
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $last.Range
$tail.Collapse(0)
$tail.InsertAfter("`r`rMay 1 2016`rDecided to rewrite the Arduino code from scratch. This time the smart way… Estimating rpms, running state etc…")

# Bold the "May 1 2016" heading paragraph (run + paragraph mark).
$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$headingPara.Range.Bold = 1

# ------------------------------------------------------------------
# 3) Split "Decided to ... scratch." off into its own run, matching
#    the target markup. A throwaway bookmark forces the run boundary
#    (adjacent same-formatted runs otherwise get merged on save).
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$split1Text = "Decided to rewrite the Arduino code from scratch."
$split1Pos = $finalPara.Range.Start + $split1Text.Length
$split1Range = $d.Range($split1Pos, $split1Pos)
$d.Bookmarks.Add("_TempSplit1", $split1Range)
$d.Bookmarks.Item("_TempSplit1").Delete()

# ------------------------------------------------------------------
# 4) Append the closing sentence (with its trailing space) as its own
#    insertion so the earlier runs don't inherit an xml:space flag
#    from this run's trailing whitespace.
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$textEnd = $finalPara.Range.End - 1
$endRange = $d.Range($textEnd, $textEnd)
$endRange.InsertAfter(" I think actually I will go with either Arduino or MK06 (when it arrives…). ")

# ------------------------------------------------------------------
# 5) Move the _GoBack bookmark from the old last paragraph to the
#    new final paragraph, right after "...etc…" and before
#    " I think actually..." (which is already its own run by now).
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$marker = " I think actually I will go with either Arduino or MK06 (when it arrives…). "
$paraText = $finalPara.Range.Text
$markerPos = $paraText.IndexOf($marker)
$bmPos = $finalPara.Range.Start + $markerPos
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
